$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "41.504.49"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +4.34%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.221.61"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "230.90"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.626"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "61.07"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.09%  "
$ws.Range("E8").Value = "  +0.01%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.403"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.99%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "58.69"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0888"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +5.36%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.549.45"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "15.66"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "21.86"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.798"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.56"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.18%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.219.69"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.54%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "41.340.96"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.96"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0899"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.81%  "
$ws.Range("E22").Value = "  +0.51%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "251.35"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +10.11%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  -2.08%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.56"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.30%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "168.29"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  +0.89%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "19.97"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("E31").Value = "  -0.98%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  -5.39%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("E40").Value = "  +28.21%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +5.32%  "
$ws.Range("E43").Value = "  -0.56%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.60"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +8.41%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0985"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +5.82%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "99.09"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("E47").Value = "  -0.06%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.465.81"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.17%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "16.60"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -1.38%  "
